$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 (Pernambuco): value and placement change
$ws.Range("C3").Value = -0.076
$ws.Range("D3").Value = "2º"

# Select cell D4, as reflected in the saved sheet view
$ws.Range("D4").Select()
